$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.191.71'
$ws.Range('E2').Value = '  +2.03%  '

$ws.Range('D3').Value = '3.123.14'
$ws.Range('E3').Value = '  +4.27%  '

$ws.Range('E4').Value = '  -0.34%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '561.15'
$c.ClearFormats()
$ws.Range('E5').Value = '  +3.94%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '139.56'
$c.ClearFormats()
$ws.Range('E6').Value = '  +4.05%  '

$ws.Range('E7').Value = '  -0.44%  '

$ws.Range('D8').Value = '3.117.45'
$ws.Range('E8').Value = '  +4.46%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.493'
$c.ClearFormats()
$ws.Range('E9').Value = '  +3.47%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '6.77'
$c.ClearFormats()
$ws.Range('E10').Value = '  +7.98%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.153'
$c.ClearFormats()
$ws.Range('E11').Value = '  +4.18%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.459'
$c.ClearFormats()
$ws.Range('E12').Value = '  +4.14%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '35.89'
$c.ClearFormats()
$ws.Range('E13').Value = '  +4.17%  '

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.0000218'
$c.ClearFormats()
$ws.Range('E14').Value = '  +3.82%  '

$ws.Range('D15').Value = '3.660.56'
$ws.Range('E15').Value = '  +4.16%  '

$ws.Range('D16').Value = '64.314.46'
$ws.Range('E16').Value = '  +1.55%  '

$ws.Range('D17').Value = '3.152.61'
$ws.Range('E17').Value = '  +3.20%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '508.19'
$c.ClearFormats()
$ws.Range('E19').Value = '  +7.84%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.73'
$c.ClearFormats()
$ws.Range('E20').Value = '  +4.89%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '13.83'
$c.ClearFormats()
$ws.Range('E21').Value = '  +4.68%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.710'
$c.ClearFormats()
$ws.Range('E22').Value = '  +6.71%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '7.32'
$c.ClearFormats()
$ws.Range('E23').Value = '  +5.22%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '12.55'
$c.ClearFormats()
$ws.Range('E24').Value = '  +4.44%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '78.48'
$c.ClearFormats()
$ws.Range('E25').Value = '  +3.15%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range('E26').Value = '  +0.46%  '

$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '2.80'
$c.ClearFormats()
$ws.Range('E27').Value = '  +5.71%  '

$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '8.52'
$c.ClearFormats()
$ws.Range('E28').Value = '  +11.29%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.07'
$c.ClearFormats()
$ws.Range('E29').Value = '  +3.27%  '

$ws.Range('E30').Value = '  -0.40%  '

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '26.39'
$c.ClearFormats()
$ws.Range('E31').Value = '  +5.40%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '2.58'
$c.ClearFormats()
$ws.Range('E32').Value = '  +1.97%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.13'
$c.ClearFormats()
$ws.Range('E33').Value = '  +3.24%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '554.49'
$c.ClearFormats()
$ws.Range('E34').Value = '  -1.15%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '55.48'
$c.ClearFormats()
$ws.Range('E35').Value = '  +8.10%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '5.99'
$c.ClearFormats()
$ws.Range('E36').Value = '  +3.77%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '5.27'
$c.ClearFormats()
$ws.Range('E37').Value = '  +0.78%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.0421'
$c.ClearFormats()
$ws.Range('E38').Value = '  +7.80%  '

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.0807'
$c.ClearFormats()
$ws.Range('E39').Value = '  +5.55%  '

$ws.Range('D40').Value = '3.086.58'
$ws.Range('E40').Value = '  +6.53%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.120'
$c.ClearFormats()
$ws.Range('E41').Value = '  +4.56%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '8.17'
$c.ClearFormats()
$ws.Range('E42').Value = '  +2.29%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '2.64'
$c.ClearFormats()
$ws.Range('E43').Value = '  -3.09%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.259'
$c.ClearFormats()
$ws.Range('E44').Value = '  +10.41%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.15'
$c.ClearFormats()
$ws.Range('E45').Value = '  +6.58%  '

$ws.Range('E46').Value = '  +0.00%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '121.54'
$c.ClearFormats()
$ws.Range('E47').Value = '  +3.49%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '24.60'
$c.ClearFormats()
$ws.Range('E48').Value = '  +2.72%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.107'
$c.ClearFormats()
$ws.Range('E49').Value = '  +2.35%  '

$ws.Range('D50').Value = '0.0₃0509'
$ws.Range('E50').Value = '  +0.58%  '

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.05'
$c.ClearFormats()
$ws.Range('E51').Value = '  +3.12%  '
